$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

$ws.Range("B2").Value = 1.334451744827586
$ws.Range("C2").Value = 1.010044170114943
$ws.Range("D2").Value = 1.160596259770115
$ws.Range("E2").Value = 1.591074064367816
$ws.Range("F2").Value = 1.008805643678161

$ws.Range("B3").Value = 2.552346855172413
$ws.Range("C3").Value = 2.306662970114942
$ws.Range("D3").Value = 1.921958181609196
$ws.Range("E3").Value = 2.376835326436782
$ws.Range("F3").Value = 1.841203227586207

$ws.Range("B4").Value = 1.098509287356322
$ws.Range("C4").Value = 1.327625956321839
$ws.Range("D4").Value = 1.349133710344827
$ws.Range("E4").Value = 1.880159213793104
$ws.Range("F4").Value = 1.251415147126436

$ws.Range("B5").Value = 3.849243278160919
$ws.Range("C5").Value = 2.198346811494253
$ws.Range("D5").Value = 2.008813149425288
$ws.Range("E5").Value = 2.688602498850574
$ws.Range("F5").Value = 1.840514965517242
